$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q4" sheet right after "总计" and before
#    "2022-Q3" (which pushes all the existing quarter sheets one slot later
#    in tab order / file order, exactly like the diff shows).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

# Use the (still existing, now shifted) "2022-Q3" sheet as a formatting
# template for the header row / leading index column so the new sheet gets
# the same bold+border style (style index 2 in the original workbook).
$q3 = $wb.Worksheets.Item("2022-Q3")

$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$q3.Cells.Item(2, 1).Copy()
$q4.Cells.Item(2, 1).PasteSpecial(-4122)

$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

# Columns B..G hold text that merely looks numeric ("8.86", "0.3083", ...) in
# the source data, so force text storage before writing the values, then
# drop back to the Normal style (keeps the "text" typing, loses the
# number-format override) to match the unstyled data cells in the diff.
$q4.Range("B2:G2").NumberFormat = "@"
$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "257010"
$q4.Cells.Item(2, 3).Value = "国联安小盘精选混合"
$q4.Cells.Item(2, 4).Value = "8.86"
$q4.Cells.Item(2, 5).Value = "74.78"
$q4.Cells.Item(2, 6).Value = "3.48"
$q4.Cells.Item(2, 7).Value = "0.3083"
$q4.Cells.Item(2, 8).Value = 10
$q4.Range("B2:G2").Style = "Normal"

# ---------------------------------------------------------------------------
# 2. Add the corresponding summary row to "总计" (sheet 1): a new row 2 for
#    2022-Q4, shifting the previous rows 2-8 down to 3-9.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")
$ws1.Rows("2:2").Insert()

# Re-use row 3's (the row that used to be row 2, now pushed down) formatting
# for the new A2 cell so it keeps the bold/border style used by the rest of
# column A, then clear the borrowed header-ish formatting that Insert()
# copied into B2:D2 so they stay plain like the other data rows.
$ws1.Cells.Item(3, 1).Copy()
$ws1.Cells.Item(2, 1).PasteSpecial(-4122)
$ws1.Range("B2:D2").ClearFormats()

$ws1.Cells.Item(2, 1).Value = 0
$ws1.Cells.Item(2, 2).Value = "2022-Q4"
$ws1.Cells.Item(2, 3).Value = 1
$ws1.Cells.Item(2, 4).Value = 0.31
